$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.229.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5211"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4397"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09279"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.644"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.126.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.878"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001155"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06709"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.318"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.08%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.271.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.503"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.126"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1048"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.661"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.710"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.205"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.890"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02615"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06752"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6967"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.346"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.12%  "

$ws.Range("E41").Value = "  +0.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2208"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6810"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.006"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.369"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.644"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000348"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.206"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.215"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "

